$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @(
    "FAPs", "Agt", "Agtr2", "FAPs",
    1, 0.3333333333333333, 0.174539, 0.523617, 0.2427616627057682, 0.2427616627057681,
    3, 1, 0.9721403333333333, 2.916421, 0.4203533802573702, 0.4203533802573701,
    0.1696764016396667, 1.527087614757, 0.1020456855152692, 0.1020456855152692
)

$row3 = @(
    "FAPs", "Agt", "Agtr2", "sCs",
    1, 0.3333333333333333, 0.174539, 0.523617, 0.2427616627057682, 0.2427616627057681,
    2, 0.6666666666666666, 1.340533666666667, 4.021601, 0.5796466197426299, 0.5796466197426299,
    0.2339754056463333, 2.105778650817, 0.140715977190499, 0.1407159771904989
)

$row4 = @(
    "sCs", "Agt", "Agtr2", "FAPs",
    3, 1, 0.5444336666666666, 1.633301, 0.7572383372942318, 0.7572383372942317,
    3, 1, 0.9721403333333333, 2.916421, 0.4203533802573702, 0.4203533802573701,
    0.5292659261912221, 4.763393335720999, 0.318307694742101, 0.3183076947421009
)

$row5 = @(
    "sCs", "Agt", "Agtr2", "sCs",
    3, 1, 0.5444336666666666, 1.633301, 0.7572383372942318, 0.7572383372942317,
    2, 0.6666666666666666, 1.340533666666667, 4.021601, 0.5796466197426299, 0.5796466197426299,
    0.7298316594334444, 6.568484934901, 0.4389306425521309, 0.4389306425521308
)

$rows = @($row2, $row3, $row4, $row5)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $data = $rows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $data.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $data[$c]
    }
}
